# Cap nhat danh sach uy quyen mau:
#  - CCCD Nguoi uy quyen o A2 duoc nhap lai nhu mot con so (mat so 0 dau)
#  - Them 2 cot tieu de moi: "Dia chi" (J1) va "Noi cap" (K1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "001090001234" (text) -> 1090001234 (number)
$ws.Range("A2").Value = 1090001234

# New header columns
$ws.Range("J1").Value = "Địa chỉ "
$ws.Range("K1").Value = "Nơi cấp "

# Match the final selection left behind by the edit
$ws.Range("B6").Select()
